$d = $word.ActiveDocument

# The target edit rewrites the second paragraph (the anagram explanation) and
# the third paragraph (which only contained "itertools.permutations") into a
# much richer set of paragraphs covering both Question 1 and a new
# Question 2 (longest palindromic substring) explanation.
#
# We replace the exact range spanning from the start of paragraph 2 through
# the end of paragraph 3 (inclusive of paragraph 3's own end-of-paragraph
# mark, since it's the last body paragraph before the sectPr) with the new
# OOXML markup for all of the replacement paragraphs. Paragraph 1
# ("Question 1: Explanation") is left completely untouched.

$p2 = $d.Paragraphs.Item(2).Range
$p3 = $d.Paragraphs.Item(3).Range

$targetRange = $d.Range($p2.Start, $p3.End)

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00671F9F" w:rsidRDefault="00682170"><w:r><w:t xml:space="preserve">In Question1, the </w:t></w:r><w:r w:rsidR="008C188F"><w:t>faction</w:t></w:r><w:r><w:t xml:space="preserve"> take</w:t></w:r><w:r w:rsidR="008C188F"><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> two strings</w:t></w:r><w:r><w:t xml:space="preserve"> (s and t)</w:t></w:r><w:r><w:t xml:space="preserve"> and find out </w:t></w:r><w:r w:rsidR="008C188F"><w:t>whether</w:t></w:r><w:r><w:t xml:space="preserve"> one</w:t></w:r><w:r><w:t xml:space="preserve"> (t)</w:t></w:r><w:r><w:t xml:space="preserve"> is anagram of the other</w:t></w:r><w:r><w:t xml:space="preserve"> (s)</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidR="00CA4F5B"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="000E2EC7"><w:t xml:space="preserve">Here, </w:t></w:r><w:r><w:t xml:space="preserve">in </w:t></w:r><w:r><w:t xml:space="preserve">the string t, all </w:t></w:r><w:r><w:t xml:space="preserve">the different </w:t></w:r><w:r><w:t xml:space="preserve">letter combinations are considered to find out anagrams. </w:t></w:r><w:r><w:t>If an anagram is found, the program returns true.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The efficiency of the program depends on how long are the given strings. However, to increase the time efficiency </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>itertools.permutations</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>() function is used</w:t></w:r><w:r><w:t xml:space="preserve"> to generate </w:t></w:r><w:r><w:t xml:space="preserve">different </w:t></w:r><w:r><w:t>letter combinations in a word</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p w:rsidR="00671F9F" w:rsidRDefault="00671F9F"><w:proofErr w:type="spellStart"/><w:r><w:t>word_combination_list</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>list(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>itertools.permutations</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>t_list</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>len</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>t_list</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)))</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">This function is helpful to eliminate unnecessary </w:t></w:r><w:r><w:t>‘</w:t></w:r><w:r><w:t>for</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:t xml:space="preserve"> loop iterations</w:t></w:r><w:r><w:t xml:space="preserve"> and complexity</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Also, it </w:t></w:r><w:r><w:t xml:space="preserve">is </w:t></w:r><w:r><w:t>written in low</w:t></w:r><w:r><w:t>er</w:t></w:r><w:r><w:t xml:space="preserve"> level languages like C, it is faster than conventional Python ‘for’ loops.</w:t></w:r><w:r><w:t xml:space="preserve"> Since I am using minimum number of variable to store data during the program running, the space efficiency also in good standing. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Question 2: Explanation</w:t></w:r></w:p><w:p><w:r><w:t>In this fu</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t>ction</w:t></w:r><w:r><w:t>, the longest palindromic substring is found in a given string.</w:t></w:r><w:r><w:t xml:space="preserve"> To achieve this, a substring is divided in to two from the center and check whether they are backwardly matching or not. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The overall efficiency of this function also depend on the length of the string. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$targetRange.InsertXML($newXml) | Out-Null
